$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")
$ws.Select()

# Re-assert the F1 header value ("View") - this mirrors the source edit that
# re-entered the cell during the merge fix.
$ws.Range("F1").Value = "View"

# Fix merge error: set the "View" flag (column F) to TRUE for rows 68-75,
# which were missing a value in this column.
$ws.Range("F68:F75").Value = $true

# Fix merge error values in rows 76-78 (column F "View" flag, plus the
# Public/Private flags in columns C/D/E that had been incorrectly merged).
$ws.Range("D76").Value = $false
$ws.Range("E76").Value = $false
$ws.Range("F76").Value = $true

$ws.Range("D77").Value = $false
$ws.Range("E77").Value = $false
$ws.Range("F77").Value = $true

$ws.Range("C78").Value = $false

# Restore the view/selection state on the sheet.
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C78").Select()
